$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 27 (shifts existing rows 27.. down by 2)
$ws.Rows("27:28").Insert()

# New row 27: "Extra" quality record dated 45238
$ws.Cells.Item(27,1).Value  = 1
$ws.Cells.Item(27,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(27,4).Value  = 45238
$ws.Cells.Item(27,5).Value  = 15
$ws.Cells.Item(27,6).Value  = 100112028
$ws.Cells.Item(27,7).Value  = "Sandia"
$ws.Cells.Item(27,8).Value  = "Sin especificar"
$ws.Cells.Item(27,9).Value  = "Extra"
$ws.Cells.Item(27,10).Value = 500
$ws.Cells.Item(27,11).Value = 450
$ws.Cells.Item(27,12).Value = 480
$ws.Cells.Item(27,13).Value = 465
$ws.Cells.Item(27,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(27,15).Value = "Perú"
$ws.Cells.Item(27,16).Value = 465
$ws.Cells.Item(27,17).Value = 1
$ws.Cells.Item(27,18).Value = "Hortaliza"

# New row 28: "Primera" quality record dated 45238
$ws.Cells.Item(28,1).Value  = 1
$ws.Cells.Item(28,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(28,4).Value  = 45238
$ws.Cells.Item(28,5).Value  = 15
$ws.Cells.Item(28,6).Value  = 100112028
$ws.Cells.Item(28,7).Value  = "Sandia"
$ws.Cells.Item(28,8).Value  = "Sin especificar"
$ws.Cells.Item(28,9).Value  = "Primera"
$ws.Cells.Item(28,10).Value = 800
$ws.Cells.Item(28,11).Value = 450
$ws.Cells.Item(28,12).Value = 480
$ws.Cells.Item(28,13).Value = 465
$ws.Cells.Item(28,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(28,15).Value = "Perú"
$ws.Cells.Item(28,16).Value = 465
$ws.Cells.Item(28,17).Value = 1
$ws.Cells.Item(28,18).Value = "Hortaliza"
